# 04_Pedidos.xlsx - adapt sheet to the full Trello "Las-Lira" flow.
# - Insert a new "Día Entrega" column right after "Fecha Entrega" (new col D),
#   shifting old D:Q -> E:R.
# - Append 4 new trailing columns: Estado Pago (S), Tipo Pedido (T),
#   Cobranza (U), Foto Enviado (V).
# - Update the "Estado" column (now R) values to the new pipeline states.
# - Fill in the new columns' data for every order row.
# - Re-apply the explicit column widths as specified by the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert the new "Día Entrega" column at D - shifts Canal..Estado right by one.
$ws.Columns("D").Insert()

# 2) New column D header + values (day-of-week label per order).
$ws.Range("D1").Value = "Día Entrega"
$ws.Range("D2").Value  = "LUNES"
$ws.Range("D3").Value  = "MARTES"
$ws.Range("D4").Value  = "MIERCOLES"
$ws.Range("D5").Value  = "MIERCOLES"
$ws.Range("D6").Value  = "JUEVES"
$ws.Range("D7").Value  = "VIERNES"
$ws.Range("D8").Value  = "JUEVES"
$ws.Range("D9").Value  = "JUEVES"
$ws.Range("D10").Value = "SABADO"

# 3) Append the 4 new trailing columns (S:V). Copy the header style from the
#    existing "Estado" header (now R1) so the new headers match formatting.
$ws.Range("R1").Copy($ws.Range("S1:V1"))

$ws.Range("S1").Value = "Estado Pago"
$ws.Range("T1").Value = "Tipo Pedido"
$ws.Range("U1").Value = "Cobranza"
$ws.Range("V1").Value = "Foto Enviado"

# 4) Update the "Estado" column (now R) to the new pipeline state names.
$ws.Range("R2").Value  = "Despachados"
$ws.Range("R3").Value  = "Despachados"
$ws.Range("R4").Value  = "Listo para Despacho"
$ws.Range("R5").Value  = "En Proceso"
$ws.Range("R6").Value  = "Entregas de Hoy"
$ws.Range("R7").Value  = "Entregas para Mañana"
$ws.Range("R8").Value  = "Entregas de Hoy"
$ws.Range("R9").Value  = "Pedidos Semana"
$ws.Range("R10").Value = "Pedidos Semana"

# 5) Estado Pago (S), Tipo Pedido (T), Cobranza (U), Foto Enviado (V) per row.
$ws.Range("S2").Value  = "Pagado"
$ws.Range("T2").Value  = "Normal"
$ws.Range("U2").Value  = "BOLETA 11248 TR. 21/10/25"

$ws.Range("S3").Value  = "Pagado"
$ws.Range("T3").Value  = "Normal"
$ws.Range("U3").Value  = "FACTURA 2345 TR. 22/10/25"
$ws.Range("V3").Value  = "arreglo_enviado_ped002.jpg"

$ws.Range("S4").Value  = "Pagado"
$ws.Range("T4").Value  = "Normal"
$ws.Range("U4").Value  = "BOLETA 11249 TR. 23/10/25"

$ws.Range("S5").Value  = "Pagado"
$ws.Range("T5").Value  = "Normal"

$ws.Range("S6").Value  = "No Pagado"
$ws.Range("T6").Value  = "Normal"

$ws.Range("S7").Value  = "Pagado"
$ws.Range("T7").Value  = "EVENTO"
$ws.Range("U7").Value  = "FACTURA 2346 TR. 25/10/25"

$ws.Range("S8").Value  = "No Pagado"
$ws.Range("T8").Value  = "Normal"

$ws.Range("S9").Value  = "Falta Boleta o Factura"
$ws.Range("T9").Value  = "Normal"

$ws.Range("S10").Value = "Pagado"
$ws.Range("T10").Value = "MANTENCIONES"

# 6) Re-apply explicit column widths for the full new layout (A:V).
#    NOTE: the interop's Range.ColumnWidth adds the standard ~5/6 character
#    padding on top of the value before it lands in the OOXML <col width=.../>
#    attribute, so we subtract that offset here to land exactly on the
#    target raw widths (10, 16, 16, 12, 10, 15, 18, 15, 25, 25, 12, 12, 18,
#    30, 20, 35, 15, 15, 18, 15, 15, 25).
$padding = 5/6
$ws.Columns("A").ColumnWidth = 10 - $padding
$ws.Columns("B").ColumnWidth = 16 - $padding
$ws.Columns("C").ColumnWidth = 16 - $padding
$ws.Columns("D").ColumnWidth = 12 - $padding
$ws.Columns("E").ColumnWidth = 10 - $padding
$ws.Columns("F").ColumnWidth = 15 - $padding
$ws.Columns("G").ColumnWidth = 18 - $padding
$ws.Columns("H").ColumnWidth = 15 - $padding
$ws.Columns("I").ColumnWidth = 25 - $padding
$ws.Columns("J").ColumnWidth = 25 - $padding
$ws.Columns("K").ColumnWidth = 12 - $padding
$ws.Columns("L").ColumnWidth = 12 - $padding
$ws.Columns("M").ColumnWidth = 18 - $padding
$ws.Columns("N").ColumnWidth = 30 - $padding
$ws.Columns("O").ColumnWidth = 20 - $padding
$ws.Columns("P").ColumnWidth = 35 - $padding
$ws.Columns("Q").ColumnWidth = 15 - $padding
$ws.Columns("R").ColumnWidth = 15 - $padding
$ws.Columns("S").ColumnWidth = 18 - $padding
$ws.Columns("T").ColumnWidth = 15 - $padding
$ws.Columns("U").ColumnWidth = 15 - $padding
$ws.Columns("V").ColumnWidth = 25 - $padding

Write-Output "Pedidos sheet updated: Dia Entrega column inserted, Estado Pago/Tipo Pedido/Cobranza/Foto Enviado appended, Estado values remapped."
